# Scheduled runner update: refresh cached market-price derived figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve-profit
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Values below were
# recomputed upstream from refreshed market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 12483.3
$ws.Range("I28").Value = 1111.6666
$ws.Range("J28").Value = 17356.857
$ws.Range("K28").Value = 1111.6666
$ws.Range("L28").Value = 17356.857
$ws.Range("M28").Value = -626.6666
$ws.Range("N28").Value = -18326.857

# Row 32
$ws.Range("H32").Value = 784
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 784
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 784
$ws.Range("N32").Value = -1436

# Row 34
$ws.Range("H34").Value = 1555.5
$ws.Range("I34").Value = 1520.5714
$ws.Range("J34").Value = 1800
$ws.Range("K34").Value = 1520.5714
$ws.Range("L34").Value = 1800
$ws.Range("M34").Value = -1317.5714
$ws.Range("N34").Value = -2206

# Row 36
$ws.Range("H36").Value = 1555.5
$ws.Range("I36").Value = 1520.5714
$ws.Range("J36").Value = 1800
$ws.Range("K36").Value = 1520.5714
$ws.Range("L36").Value = 1800
$ws.Range("M36").Value = -805.5714
$ws.Range("N36").Value = -3230

# Row 64
$ws.Range("H64").Value = 7118
$ws.Range("I64").Value = 5899.8335
$ws.Range("J64").Value = 8579.799999999999
$ws.Range("K64").Value = 5899.8335
$ws.Range("L64").Value = 8579.799999999999
$ws.Range("M64").Value = -5651.8335
$ws.Range("N64").Value = -9075.799999999999

# Row 67
$ws.Range("H67").Value = 7118
$ws.Range("I67").Value = 5899.8335
$ws.Range("J67").Value = 8579.799999999999
$ws.Range("K67").Value = 5899.8335
$ws.Range("L67").Value = 8579.799999999999
$ws.Range("M67").Value = -5041.8335
$ws.Range("N67").Value = -10295.8

# Row 112
$ws.Range("H112").Value = 940.8
$ws.Range("I112").Value = 994.5
$ws.Range("J112").Value = 726
$ws.Range("K112").Value = 2983.5
$ws.Range("L112").Value = 2178
$ws.Range("M112").Value = -1875.5
$ws.Range("N112").Value = -4394

# Row 132
$ws.Range("H132").Value = 16493.75
$ws.Range("I132").Value = 16284.615
$ws.Range("J132").Value = 17400
$ws.Range("K132").Value = 48853.845
$ws.Range("L132").Value = 52200
$ws.Range("M132").Value = -46323.845
$ws.Range("N132").Value = -57260

# Row 138
$ws.Range("H138").Value = 1973.3103
$ws.Range("I138").Value = 1027.8422
$ws.Range("J138").Value = 3769.7
$ws.Range("K138").Value = 3083.5266
$ws.Range("L138").Value = 11309.1
$ws.Range("M138").Value = 2056.4734
$ws.Range("N138").Value = -21589.1

$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 96
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 96
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 191

# Row 61
$ws.Range("H61").Value = 7000
$ws.Range("I61").Value = 7000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6788

# Row 74
$ws.Range("H74").Value = 8559.799999999999
$ws.Range("I74").Value = 8000
$ws.Range("J74").Value = 8933
$ws.Range("K74").Value = 8000
$ws.Range("L74").Value = 8933
$ws.Range("M74").Value = -7126
$ws.Range("N74").Value = -10681

# Row 77
$ws.Range("H77").Value = 8559.799999999999
$ws.Range("I77").Value = 8000
$ws.Range("J77").Value = 8933
$ws.Range("K77").Value = 40000
$ws.Range("L77").Value = 44665
$ws.Range("M77").Value = -35632
$ws.Range("N77").Value = -53401

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()

# Row 136
$ws.Range("H136").Value = 7000
$ws.Range("I136").Value = 7000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 21000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -18450

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 178.88235
$ws.Range("I80").Value = 127.166664
$ws.Range("J80").Value = 207.09091
$ws.Range("K80").Value = 127.166664
$ws.Range("L80").Value = 207.09091
$ws.Range("M80").Value = 870.833336
$ws.Range("N80").Value = -2203.09091

# Row 83
$ws.Range("H83").Value = 178.88235
$ws.Range("I83").Value = 127.166664
$ws.Range("J83").Value = 207.09091
$ws.Range("K83").Value = 635.83332
$ws.Range("L83").Value = 1035.45455
$ws.Range("M83").Value = 4356.16668
$ws.Range("N83").Value = -11019.45455

# Row 86
$ws.Range("H86").Value = 4721.4165
$ws.Range("I86").Value = 2174.75
$ws.Range("J86").Value = 5994.75
$ws.Range("K86").Value = 2174.75
$ws.Range("L86").Value = 5994.75
$ws.Range("M86").Value = -1051.75
$ws.Range("N86").Value = -8240.75

# Row 89
$ws.Range("H89").Value = 4721.4165
$ws.Range("I89").Value = 2174.75
$ws.Range("J89").Value = 5994.75
$ws.Range("K89").Value = 10873.75
$ws.Range("L89").Value = 29973.75
$ws.Range("M89").Value = -5257.75
$ws.Range("N89").Value = -41205.75

# Row 107
$ws.Range("H107").Value = 4018.56
$ws.Range("I107").Value = 1283.2142
$ws.Range("J107").Value = 7499.909
$ws.Range("K107").Value = 1283.2142
$ws.Range("L107").Value = 7499.909
$ws.Range("M107").Value = 636.7858000000001
$ws.Range("N107").Value = -11339.909

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3466.5881
$ws.Range("I58").Value = 866.5
$ws.Range("J58").Value = 4884.8184
$ws.Range("K58").Value = 866.5
$ws.Range("L58").Value = 4884.8184
$ws.Range("M58").Value = -663.5
$ws.Range("N58").Value = -5290.8184

# Row 99
$ws.Range("H99").Value = 3823.25
$ws.Range("I99").Value = 3342.2856
$ws.Range("J99").Value = 4496.6
$ws.Range("K99").Value = 3342.2856
$ws.Range("L99").Value = 4496.6
$ws.Range("M99").Value = -1844.2856
$ws.Range("N99").Value = -7492.6

# Row 122
$ws.Range("H122").Value = 1697.8334
$ws.Range("I122").Value = 1583
$ws.Range("J122").Value = 1812.6666
$ws.Range("K122").Value = 4749
$ws.Range("L122").Value = 5437.9998
$ws.Range("M122").Value = -2299
$ws.Range("N122").Value = -10337.9998

# Row 126
$ws.Range("H126").Value = 3823.25
$ws.Range("I126").Value = 3342.2856
$ws.Range("J126").Value = 4496.6
$ws.Range("K126").Value = 10026.8568
$ws.Range("L126").Value = 13489.8
$ws.Range("M126").Value = -7556.856800000001
$ws.Range("N126").Value = -18429.8

# Row 136
$ws.Range("H136").Value = 3466.5881
$ws.Range("I136").Value = 866.5
$ws.Range("J136").Value = 4884.8184
$ws.Range("K136").Value = 2599.5
$ws.Range("L136").Value = 14654.4552
$ws.Range("M136").Value = -49.5
$ws.Range("N136").Value = -19754.4552

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 329.14285
$ws.Range("I12").Value = 158.33333
$ws.Range("J12").Value = 375.72726
$ws.Range("K12").Value = 474.99999
$ws.Range("L12").Value = 1127.18178
$ws.Range("M12").Value = -301.99999
$ws.Range("N12").Value = -1473.18178

# Row 42
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

# Row 61
$ws.Range("H61").Value = 75
$ws.Range("I61").Value = 100
$ws.Range("J61").Value = 50
$ws.Range("K61").Value = 300
$ws.Range("L61").Value = 150
$ws.Range("M61").Value = -85
$ws.Range("N61").Value = -580

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1399.6666
$ws.Range("I97").Value = 1224.5
$ws.Range("J97").Value = 1750
$ws.Range("K97").Value = 1224.5
$ws.Range("L97").Value = 1750
$ws.Range("M97").Value = -728.5
$ws.Range("N97").Value = -2742

# Row 113
$ws.Range("H113").Value = 8000.5
$ws.Range("I113").Value = 5335
$ws.Range("J113").Value = 9999.625
$ws.Range("K113").Value = 5335
$ws.Range("L113").Value = 9999.625
$ws.Range("M113").Value = -3165
$ws.Range("N113").Value = -14339.625

# Row 122
$ws.Range("H122").Value = 334542.6
$ws.Range("I122").Value = 455785.2
$ws.Range("J122").Value = 1125.5
$ws.Range("K122").Value = 1367355.6
$ws.Range("L122").Value = 3376.5
$ws.Range("M122").Value = -1364905.6
$ws.Range("N122").Value = -8276.5

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 999999
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 999999
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 999999
$ws.Range("N20").Value = -1000451

# Row 22
$ws.Range("H22").Value = 599
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 599
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 599
$ws.Range("N22").Value = -1189

# Row 27
$ws.Range("H27").Value = 599
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 599
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 599
$ws.Range("N27").Value = -813

# Row 40
$ws.Range("H40").Value = 4599
$ws.Range("I40").Value = 2716.4443
$ws.Range("J40").Value = 7987.6
$ws.Range("K40").Value = 2716.4443
$ws.Range("L40").Value = 7987.6
$ws.Range("M40").Value = -2580.4443
$ws.Range("N40").Value = -8259.6

# Row 55
$ws.Range("H55").Value = 1077.1818
$ws.Range("I55").Value = 1591.5
$ws.Range("J55").Value = 460
$ws.Range("K55").Value = 1591.5
$ws.Range("L55").Value = 460
$ws.Range("M55").Value = -1418.5
$ws.Range("N55").Value = -806

# Row 136
$ws.Range("H136").Value = 3457
$ws.Range("I136").Value = 3457
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10371
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7821
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 3800.4
$ws.Range("I100").Value = 3800.4
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 7600.8
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -7059.8

# Row 113
$ws.Range("H113").Value = 543.1667
$ws.Range("I113").Value = 631.6667
$ws.Range("J113").Value = 366.16666
$ws.Range("K113").Value = 1895.0001
$ws.Range("L113").Value = 1098.49998
$ws.Range("M113").Value = 274.9999
$ws.Range("N113").Value = -5438.499980000001

# Row 136
$ws.Range("H136").Value = 8312.5
$ws.Range("I136").Value = 4400
$ws.Range("J136").Value = 9616.666999999999
$ws.Range("K136").Value = 13200
$ws.Range("L136").Value = 28850.001
$ws.Range("M136").Value = -10650
$ws.Range("N136").Value = -33950.001
